# Auto-generated edit script applying scheduled-runner profit recalculations
# to the Sagittarius_Profits workbook (per-sheet Leve profit tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 6 (item id row)
$ws.Range("H6").Value = 40.764706
$ws.Range("I6").Value = 40.1875
$ws.Range("K6").Value = 120.5625
$ws.Range("M6").Value = -8.5625

# ALC!row 29 (item id row)
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# ALC!row 52 (item id row)
$ws.Range("H52").Value = 190
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# ALC!row 62 (item id row)
$ws.Range("H62").Value = 8590.923000000001
$ws.Range("I62").Value = 8468.200000000001
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 8468.200000000001
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -7844.200000000001
$ws.Range("N62").Value = -10248

# ALC!row 64 (item id row)
$ws.Range("H64").Value = 5000
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# ALC!row 65 (item id row)
$ws.Range("H65").Value = 8590.923000000001
$ws.Range("I65").Value = 8468.200000000001
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 42341
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -39221
$ws.Range("N65").Value = -51240

# ALC!row 67 (item id row)
$ws.Range("H67").Value = 5000
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# ALC!row 137 (item id row)
$ws.Range("H137").Value = 6103.8
$ws.Range("I137").Value = 5504.75
$ws.Range("K137").Value = 16514.25
$ws.Range("M137").Value = -13964.25

# ALC!row 138 (item id row)
$ws.Range("H138").Value = 6859.469
$ws.Range("J138").Value = 7179.909
$ws.Range("L138").Value = 21539.727
$ws.Range("N138").Value = -31819.727

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 28 (item id row)
$ws.Range("H28").Value = 13390
$ws.Range("I28").Value = 13390
$ws.Range("K28").Value = 13390
$ws.Range("M28").Value = -13198

# ARM!row 61 (item id row)
$ws.Range("H61").Value = 6214.143
$ws.Range("J61").Value = 7590.6
$ws.Range("L61").Value = 7590.6
$ws.Range("N61").Value = -8014.6

# ARM!row 63 (item id row)
$ws.Range("H63").Value = 2883.5
$ws.Range("I63").Value = 2883.5
$ws.Range("K63").Value = 2883.5
$ws.Range("M63").Value = -2197.5

# ARM!row 66 (item id row)
$ws.Range("H66").Value = 2883.5
$ws.Range("I66").Value = 2883.5
$ws.Range("K66").Value = 14417.5
$ws.Range("M66").Value = -10985.5

# ARM!row 97 (item id row)
$ws.Range("H97").Value = 1468.6923
$ws.Range("I97").Value = 1326.6364
$ws.Range("K97").Value = 1326.6364
$ws.Range("M97").Value = -830.6364000000001

# ARM!row 99 (item id row)
$ws.Range("H99").Value = 13390
$ws.Range("I99").Value = 13390
$ws.Range("K99").Value = 13390
$ws.Range("M99").Value = -10395

# ARM!row 113 (item id row)
$ws.Range("H113").Value = 150000
$ws.Range("J113").Value = 150000
$ws.Range("L113").Value = 150000
$ws.Range("N113").Value = -158678

# ARM!row 122 (item id row)
$ws.Range("H122").Value = 2914.8333
$ws.Range("I122").Value = 2897.8
$ws.Range("K122").Value = 8693.400000000001
$ws.Range("M122").Value = -6243.400000000001

# ARM!row 136 (item id row)
$ws.Range("H136").Value = 6214.143
$ws.Range("J136").Value = 7590.6
$ws.Range("L136").Value = 22771.8
$ws.Range("N136").Value = -27871.8

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 94 (item id row)
$ws.Range("H94").Value = 2050
$ws.Range("J94").Value = 1100
$ws.Range("L94").Value = 1100
$ws.Range("N94").Value = -2002

# BSM!row 99 (item id row)
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# BSM!row 107 (item id row)
$ws.Range("H107").Value = 1999
$ws.Range("I107").Value = 1999
$ws.Range("K107").Value = 1999
$ws.Range("M107").Value = -79

# BSM!row 134 (item id row)
$ws.Range("H134").Value = 4218.3076
$ws.Range("I134").Value = 4113.778
$ws.Range("K134").Value = 12341.334
$ws.Range("M134").Value = -9806.334000000001

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 58 (item id row)
$ws.Range("H58").Value = 1749.5
$ws.Range("I58").Value = 1749.5
$ws.Range("K58").Value = 1749.5
$ws.Range("M58").Value = -1546.5

# CRP!row 134 (item id row)
$ws.Range("H134").Value = 4148.4165
$ws.Range("I134").Value = 4228.1
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 12684.3
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -10149.3
$ws.Range("N134").Value = -16320

# CRP!row 136 (item id row)
$ws.Range("H136").Value = 1749.5
$ws.Range("I136").Value = 1749.5
$ws.Range("K136").Value = 5248.5
$ws.Range("M136").Value = -2698.5

# CRP!row 141 (item id row)
$ws.Range("H141").Value = 497508.94
$ws.Range("J141").Value = 497508.94
$ws.Range("L141").Value = 497508.94
$ws.Range("N141").Value = -507868.94

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 131 (item id row)
$ws.Range("H131").Value = 2958.8
$ws.Range("I131").Value = 1990
$ws.Range("K131").Value = 5970
$ws.Range("M131").Value = -930

# CUL!row 132 (item id row)
$ws.Range("H132").Value = 14316.667

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 70 (item id row)
$ws.Range("H70").Value = 5966.4375
$ws.Range("I70").Value = 5888.6
$ws.Range("K70").Value = 5888.6
$ws.Range("M70").Value = -5618.6

# GSM!row 73 (item id row)
$ws.Range("H73").Value = 5966.4375
$ws.Range("I73").Value = 5888.6
$ws.Range("K73").Value = 5888.6
$ws.Range("M73").Value = -4952.6

# GSM!row 80 (item id row)
$ws.Range("H80").Value = 4995.6665
$ws.Range("J80").Value = 5096.364
$ws.Range("L80").Value = 5096.364
$ws.Range("N80").Value = -7092.364

# GSM!row 83 (item id row)
$ws.Range("H83").Value = 4995.6665
$ws.Range("J83").Value = 5096.364
$ws.Range("L83").Value = 25481.82
$ws.Range("N83").Value = -35465.82

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 46 (item id row)
$ws.Range("H46").Value = 2667.0833
$ws.Range("I46").Value = 2519.375
$ws.Range("J46").Value = 2962.5
$ws.Range("K46").Value = 2519.375
$ws.Range("L46").Value = 2962.5
$ws.Range("M46").Value = -2331.375
$ws.Range("N46").Value = -3338.5

# LTW!row 68 (item id row)
$ws.Range("H68").Value = 2282.6667
$ws.Range("I68").Value = 3200
$ws.Range("J68").Value = 2099.2
$ws.Range("K68").Value = 3200
$ws.Range("L68").Value = 2099.2
$ws.Range("M68").Value = -2451
$ws.Range("N68").Value = -3597.2

# LTW!row 71 (item id row)
$ws.Range("H71").Value = 2282.6667
$ws.Range("I71").Value = 3200
$ws.Range("J71").Value = 2099.2
$ws.Range("K71").Value = 16000
$ws.Range("L71").Value = 10496
$ws.Range("M71").Value = -12256
$ws.Range("N71").Value = -17984

# LTW!row 132 (item id row)
$ws.Range("H132").Value = 5316.825
$ws.Range("I132").Value = 4404.4614
$ws.Range("J132").Value = 7011.2144
$ws.Range("K132").Value = 13213.3842
$ws.Range("L132").Value = 21033.6432
$ws.Range("M132").Value = -10683.3842
$ws.Range("N132").Value = -26093.6432

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 132 (item id row)
$ws.Range("H132").Value = 1866.2
$ws.Range("I132").Value = 1866.2
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5598.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3068.6
$ws.Range("N132").ClearContents()

# WVR!row 135 (item id row)
$ws.Range("H135").Value = 217147
$ws.Range("J135").Value = 217147
$ws.Range("L135").Value = 217147
$ws.Range("N135").Value = -227287

